$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''22.550.24'
$ws.Range("E2").Value = '  +0.32%  '

$ws.Range("D3").Value = '''1.577.92'
$ws.Range("E3").Value = '  +0.33%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  +0.02%  '

$ws.Range("D6").Value = '''288.75'
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("D7").Value = '''0.3688'
$ws.Range("E7").Value = '  -0.55%  '

$ws.Range("D8").Value = '''48.60'
$ws.Range("E8").Value = '  -2.73%  '

$ws.Range("D9").Value = '''0.3351'
$ws.Range("E9").Value = '  -0.84%  '

$ws.Range("D10").Value = '''1.147'
$ws.Range("E10").Value = '  +0.24%  '

$ws.Range("D11").Value = '''0.07479'
$ws.Range("E11").Value = '  -0.99%  '

$ws.Range("D13").Value = '''21.02'
$ws.Range("E13").Value = '  -1.11%  '

$ws.Range("D14").Value = '''6.010'
$ws.Range("E14").Value = '  -0.44%  '

$ws.Range("D15").Value = '''6.964'
$ws.Range("E15").Value = '  -0.05%  '

$ws.Range("D16").Value = '''1.580.35'
$ws.Range("E16").Value = '  +0.61%  '

$ws.Range("D17").Value = '''0.00001119'
$ws.Range("E17").Value = '  -0.21%  '

$ws.Range("D18").Value = '''88.78'
$ws.Range("E18").Value = '  -2.18%  '

$ws.Range("D19").Value = '''0.06772'
$ws.Range("E19").Value = '  +0.24%  '

$ws.Range("D20").Value = '''6.437'
$ws.Range("E20").Value = '  +2.28%  '

$ws.Range("E21").Value = '  +0.06%  '

$ws.Range("D22").Value = '''16.60'
$ws.Range("E22").Value = '  +1.14%  '

$ws.Range("E23").Value = '  -0.07%  '

$ws.Range("D24").Value = '''22.543.69'
$ws.Range("E24").Value = '  +0.27%  '

$ws.Range("D25").Value = '''2.399'
$ws.Range("E25").Value = '  +1.08%  '

$ws.Range("D26").Value = '''2.616'
$ws.Range("E26").Value = '  -0.36%  '

$ws.Range("D27").Value = '''152.88'
$ws.Range("E27").Value = '  +2.37%  '

$ws.Range("D28").Value = '''19.72'
$ws.Range("E28").Value = '  -1.59%  '

$ws.Range("D29").Value = '''5.019'
$ws.Range("E29").Value = '  -0.70%  '

$ws.Range("D30").Value = '''124.42'

$ws.Range("D31").Value = '''1.758.35'
$ws.Range("E31").Value = '  +0.66%  '

$ws.Range("D32").Value = '''1.072'
$ws.Range("E32").Value = '  -0.60%  '

$ws.Range("D33").Value = '''6.204'
$ws.Range("E33").Value = '  -0.74%  '

$ws.Range("D34").Value = '''2.006'
$ws.Range("E34").Value = '  -0.52%  '

$ws.Range("D35").Value = '''9.744'
$ws.Range("E35").Value = '  -0.64%  '

$ws.Range("D36").Value = '''0.08322'
$ws.Range("E36").Value = '  -0.57%  '

$ws.Range("D37").Value = '''0.02463'
$ws.Range("E37").Value = '  -0.93%  '

$ws.Range("D38").Value = '''0.2273'
$ws.Range("E38").Value = '  -1.36%  '

$ws.Range("D39").Value = '''5.457'
$ws.Range("E39").Value = '  +0.05%  '

$ws.Range("B40").Value = 'TrustWalletToken'
$ws.Range("C40").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D40").Value = '''1.309'
$ws.Range("E40").Value = '  -2.84%  '

$ws.Range("B41").Value = 'Hedera'
$ws.Range("C41").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D41").Value = '''0.06386'
$ws.Range("E41").Value = '  -2.64%  '

$ws.Range("D42").Value = '''11.42'
$ws.Range("E42").Value = '  +0.44%  '

$ws.Range("D43").Value = '''0.6358'
$ws.Range("E43").Value = '  +2.13%  '

$ws.Range("D45").Value = '''14.01'
$ws.Range("E45").Value = '  -0.28%  '

$ws.Range("D46").Value = '''0.6192'
$ws.Range("E46").Value = '  +5.68%  '

$ws.Range("E47").Value = '  -1.02%  '

$ws.Range("D48").Value = '''2.064'
$ws.Range("E48").Value = '  -0.45%  '

$ws.Range("D49").Value = '''125.29'
$ws.Range("E49").Value = '  -2.86%  '

$ws.Range("D50").Value = '''1.220'
$ws.Range("E50").Value = '  -0.04%  '

$ws.Range("E51").Value = '  -0.71%  '
